$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.412.08"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.309.44"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.72"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.66"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.32%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.310.91"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.468"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.87"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("E11").Value = "  -3.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.408"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.879.82"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.92%  "
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.90"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.312.12"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000165"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.441.56"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.18"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.37"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.68"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "374.01"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.51"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.539"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.445.27"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000102"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -7.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.172"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.19"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.36%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.61"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.04"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.60"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("E35").Value = "  -7.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.16"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.25%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.08"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.53"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.72"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("B40").Value = "RenzoRestakedETH"
$ws.Range("C40").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.341.29"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.06%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.86"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -13.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0733"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.06"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.753"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.17"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.15%  "
$ws.Range("E46").Value = "  -4.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.11"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.374.37"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.35%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.47"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -6.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.34"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.03%  "
